$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.235148
$ws.Range("H2").Value = 21.705444
$ws.Range("I2").Value = 0.9254344869740032
$ws.Range("J2").Value = 0.9254344869740032
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 30.345835
$ws.Range("N2").Value = 91.03750500000001
$ws.Range("O2").Value = 0.8527782452855476
$ws.Range("P2").Value = 0.8527782452855475
$ws.Range("Q2").Value = 219.55660740858
$ws.Range("R2").Value = 1976.00946667722
$ws.Range("S2").Value = 0.7891903979284214
$ws.Range("T2").Value = 0.7891903979284213

# row3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.235148
$ws.Range("H3").Value = 21.705444
$ws.Range("I3").Value = 0.9254344869740032
$ws.Range("J3").Value = 0.9254344869740032
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.224107666666666
$ws.Range("N3").Value = 9.672322999999999
$ws.Range("O3").Value = 0.09060383010029813
$ws.Range("P3").Value = 0.09060383010029811
$ws.Range("Q3").Value = 23.326896136268
$ws.Range("R3").Value = 209.942065226412
$ws.Range("S3").Value = 0.08384790902674914
$ws.Range("T3").Value = 0.08384790902674913

# row4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.235148
$ws.Range("H4").Value = 21.705444
$ws.Range("I4").Value = 0.9254344869740032
$ws.Range("J4").Value = 0.9254344869740032
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.014730333333333
$ws.Range("N4").Value = 6.044191
$ws.Range("O4").Value = 0.05661792461415433
$ws.Range("P4").Value = 0.05661792461415433
$ws.Range("Q4").Value = 14.576872141756
$ws.Range("R4").Value = 131.191849275804
$ws.Range("S4").Value = 0.0523961800188327
$ws.Range("T4").Value = 0.0523961800188327

# row5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.3016356666666667
$ws.Range("H5").Value = 0.9049070000000001
$ws.Range("I5").Value = 0.03858166390441884
$ws.Range("J5").Value = 0.03858166390441884
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.345835
$ws.Range("N5").Value = 91.03750500000001
$ws.Range("O5").Value = 0.8527782452855476
$ws.Range("P5").Value = 0.8527782452855475
$ws.Range("Q5").Value = 9.153386170781669
$ws.Range("R5").Value = 82.38047553703503
$ws.Range("S5").Value = 0.03290160364460705
$ws.Range("T5").Value = 0.03290160364460704

# row6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.3016356666666667
$ws.Range("H6").Value = 0.9049070000000001
$ws.Range("I6").Value = 0.03858166390441884
$ws.Range("J6").Value = 0.03858166390441884
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.224107666666666
$ws.Range("N6").Value = 9.672322999999999
$ws.Range("O6").Value = 0.09060383010029813
$ws.Range("P6").Value = 0.09060383010029811
$ws.Range("Q6").Value = 0.9725058654401111
$ws.Range("R6").Value = 8.752552788960999
$ws.Range("S6").Value = 0.00349564652138277
$ws.Range("T6").Value = 0.003495646521382769

# row7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.3016356666666667
$ws.Range("H7").Value = 0.9049070000000001
$ws.Range("I7").Value = 0.03858166390441884
$ws.Range("J7").Value = 0.03858166390441884
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.014730333333333
$ws.Range("N7").Value = 6.044191
$ws.Range("O7").Value = 0.05661792461415433
$ws.Range("P7").Value = 0.05661792461415433
$ws.Range("Q7").Value = 0.6077145272485556
$ws.Range("R7").Value = 5.469430745237
$ws.Range("S7").Value = 0.002184413738429025
$ws.Range("T7").Value = 0.002184413738429025

# row8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.2813256666666666
$ws.Range("H8").Value = 0.843977
$ws.Range("I8").Value = 0.0359838491215779
$ws.Range("J8").Value = 0.0359838491215779
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 30.345835
$ws.Range("N8").Value = 91.03750500000001
$ws.Range("O8").Value = 0.8527782452855476
$ws.Range("P8").Value = 0.8527782452855475
$ws.Range("Q8").Value = 8.537062261931666
$ws.Range("R8").Value = 76.83356035738501
$ws.Range("S8").Value = 0.03068624371251909
$ws.Range("T8").Value = 0.03068624371251909

# row9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.2813256666666666
$ws.Range("H9").Value = 0.843977
$ws.Range("I9").Value = 0.0359838491215779
$ws.Range("J9").Value = 0.0359838491215779
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.224107666666666
$ws.Range("N9").Value = 9.672322999999999
$ws.Range("O9").Value = 0.09060383010029813
$ws.Range("P9").Value = 0.09060383010029811
$ws.Range("Q9").Value = 0.907024238730111
$ws.Range("R9").Value = 8.163218148570998
$ws.Range("S9").Value = 0.003260274552166206
$ws.Range("T9").Value = 0.003260274552166206

# row10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.2813256666666666
$ws.Range("H10").Value = 0.843977
$ws.Range("I10").Value = 0.0359838491215779
$ws.Range("J10").Value = 0.0359838491215779
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.014730333333333
$ws.Range("N10").Value = 6.044191
$ws.Range("O10").Value = 0.05661792461415433
$ws.Range("P10").Value = 0.05661792461415433
$ws.Range("Q10").Value = 0.5667953541785554
$ws.Range("R10").Value = 5.101158187606999
$ws.Range("S10").Value = 0.002037330856892601
$ws.Range("T10").Value = 0.002037330856892601
